$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "91.386.94"
$ws.Range("E2").Value = "  +2.61%  "
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "3.123.98"
$ws.Range("E3").Value = "  +0.81%  "
$ws.Range("E4").Value = "  +0.04%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "218.90"
$ws.Range("E5").Value = "  +2.79%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "622.86"
$ws.Range("E6").Value = "  -0.06%  "
$ws.Range("E7").Value = "  +26.10%  "
$ws.Range("E8").Value = "  -0.30%  "
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "1.00"
$ws.Range("E9").Value = "  +0.06%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "3.121.39"
$ws.Range("E10").Value = "  +0.89%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.741"
$ws.Range("E11").Value = "  +23.66%  "
$ws.Range("E12").Value = "  +6.59%  "
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "0.0000254"
$ws.Range("E13").Value = "  +4.82%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "34.84"
$ws.Range("E14").Value = "  +7.69%  "
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "5.52"
$ws.Range("E15").Value = "  +4.10%  "
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "91.125.38"
$ws.Range("E16").Value = "  +2.93%  "
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "3.694.16"
$ws.Range("E17").Value = "  +0.72%  "
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "3.91"
$ws.Range("E18").Value = "  +15.59%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "3.098.72"
$ws.Range("E19").Value = "  +0.56%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "0.0000219"
$ws.Range("E20").Value = "  +3.42%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "14.14"
$ws.Range("E21").Value = "  +5.03%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "440.28"
$ws.Range("E22").Value = "  +3.82%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "8.87"
$ws.Range("E23").Value = "  +6.94%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "5.20"
$ws.Range("E24").Value = "  +5.25%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "6.24"
$ws.Range("E25").Value = "  +10.38%  "
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "88.55"
$ws.Range("E26").Value = "  +7.47%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "12.45"
$ws.Range("E27").Value = "  +3.95%  "
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "3.280.96"
$ws.Range("E28").Value = "  +1.52%  "
$ws.Range("E29").Value = "  -0.02%  "
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "0.167"
$ws.Range("E30").Value = "  -2.00%  "
$ws.Range("E31").Value = "  +13.11%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "527.22"
$ws.Range("E32").Value = "  +3.01%  "
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "0.898"
$ws.Range("E33").Value = "  -16.64%  "
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "3.76"
$ws.Range("E34").Value = "  +1.95%  "
$ws.Range("E35").Value = "  +13.23%  "
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "7.11"
$ws.Range("E36").Value = "  +4.68%  "
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "24.01"
$ws.Range("E37").Value = "  +7.44%  "
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "1.31"
$ws.Range("E38").Value = "  +3.74%  "
$ws.Range("E39").Value = "  +3.42%  "
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "0.0870"
$ws.Range("E40").Value = "  +25.09%  "
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "22.29"
$ws.Range("E41").Value = "  +0.05%  "
$ws.Range("E42").Value = "  +0.00%  "
$ws.Range("E43").Value = "  +18.24%  "
$ws.Range("E44").Value = "  +9.13%  "
$ws.Range("E45").Value = "  +5.94%  "
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "148.84"
$ws.Range("E47").Value = "  +1.75%  "
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "44.03"
$ws.Range("E48").Value = "  +1.63%  "
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "1.31"
$ws.Range("E49").Value = "  +7.26%  "
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "168.86"
$ws.Range("E50").Value = "  +3.49%  "
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "4.26"
$ws.Range("E51").Value = "  +7.67%  "
